$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'AgqO9cblRGPfP'
$ws.Range("B2").Value = 'Bsq7aTk3WfJVuI'
$ws.Range("A3").Value = 'yB8y2eX4O2B9i3QtP'
$ws.Range("B3").Value = '0iIB0vgKsd5Uxb'
$ws.Range("A4").Value = 'vnbJytzGhqNDj'
$ws.Range("B4").Value = 'rKxPzTWlZdoIax'
$ws.Range("A5").Value = 'eDSVmNZleIZ'
$ws.Range("B5").Value = 'eewBTLbrOdhwAn'
$ws.Range("A6").Value = 'sCTffQ9abI9uyv'
$ws.Range("B6").Value = 'Hu8n2dReuN1W9A'
$ws.Range("A7").Value = 'CjNz9WPKDM3slT'
$ws.Range("B7").Value = '6AMzxTs1dOJ5RE'
$ws.Range("A8").Value = 'RJTxKKNj'
$ws.Range("B8").Value = 'TFarAhEFenvDLC'
$ws.Range("A9").Value = 'XHZ1fSzVomfTu0G8'
$ws.Range("B9").Value = 'ukab36PJ6oBEZg'
$ws.Range("A10").Value = 'aNxBoN23Tjf'
$ws.Range("B10").Value = 'CY6Fj98hfEjNWj'
$ws.Range("A11").Value = '2Ktf6aBQzOJ'
$ws.Range("B11").Value = 'Co0zPpgQt18Zos'
$ws.Range("A12").Value = 'MuwCUV9brr'
$ws.Range("B12").Value = 'an0QXwEdV0nl5g'
$ws.Range("A13").Value = 'KVhMhGZC5J7LB'
$ws.Range("B13").Value = 'wzvCPfQ5NYxV0s'
$ws.Range("A14").Value = '9uycMVsKMldi'
$ws.Range("B14").Value = 'xc8u97rLz2YJF8'
$ws.Range("A15").Value = 'ruLDtOxVAdAeX46I'
$ws.Range("B15").Value = 'XHuj5xE23l7mYp'
$ws.Range("A16").Value = 'QyGmJibYcdfq'
$ws.Range("B16").Value = 'X9hnRe9hsmRjEG'
$ws.Range("A17").Value = 'U5ZqMJn230'
$ws.Range("B17").Value = 'YqqibHjGZicmw1'
$ws.Range("A18").Value = 'QhK0CUH1V04VT'
$ws.Range("B18").Value = '2i2kXwwI6jUPWw'
$ws.Range("A19").Value = 'Tb4haWmF'
$ws.Range("B19").Value = 'K0eC7S4JNknY8J'
$ws.Range("A20").Value = 'w1OppE1NcP'
$ws.Range("B20").Value = 'JvnuFXDNYjKPRM'
$ws.Range("A21").Value = 'aoDrfxtjiYc15'
$ws.Range("B21").Value = 'FYq9LTLyYLzDwh'
$ws.Range("A22").Value = '4VPM4BIgPjCmkUyUA'
$ws.Range("B22").Value = '3gthMNwXQYTQhq'
$ws.Range("A23").Value = 'Kmnab49t0K6J'
$ws.Range("B23").Value = '2UeON4OjdtG38m'
$ws.Range("A24").Value = '1kA5hKs3yFu6'
$ws.Range("B24").Value = 'styRoh1ShbqIkX'
$ws.Range("A25").Value = 'dfc5tLnL5l'
$ws.Range("B25").Value = '32X1rv92UoW5PV'
$ws.Range("A26").Value = 'alN2QAtmdyV'
$ws.Range("B26").Value = 'EUqwVXZ7ni9CDH'
$ws.Range("A27").Value = 'SqqWoInj8P'
$ws.Range("B27").Value = 'COUBuaMliNoHaT'
$ws.Range("A28").Value = 'bCcb6SHybod8imYz'
$ws.Range("B28").Value = 'UqigEDaswPoJWf'
$ws.Range("A29").Value = 'jjmI0Lo4ehAe'
$ws.Range("B29").Value = 'Sa1CiHhAOVD5ev'
$ws.Range("A30").Value = '9wSmdK7CAa2CQw9A'
$ws.Range("B30").Value = 'n5X928UHxGVTmk'
$ws.Range("A31").Value = 'sAEgGqW61qIqhMJ'
$ws.Range("B31").Value = '13cGvrTnAUVSfJ'
